{"js": "// The document contains three paragraphs with a transcribed\n// \"<id>p074r_N</id>\" tag. Each of these is currently split across three\n// separate runs: \"<id>\" (Courier New, dark-yellow), \"p074r_N\" (plain,\n// black) and \"</id>\" (Courier New, dark-yellow). The edit merges each\n// triplet back into a single run, e.g. \"<id>p074r_1</id>\", carrying the\n// formatting of the first (\"<id>\") run - matching the newly downloaded\n// transcription-line (tl) formatting.\nconst body = context.document.body;\n\nfor (const n of [1, 2, 3]) {\n  const fullText = `<id>p074r_${n}</id>`;\n\n  // search() matches across run boundaries, returning a single Range\n  // spanning the old \"<id>\" + \"p074r_N\" + \"</id>\" runs.\n  const results = body.search(fullText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const found of results.items) {\n    // Replacing the whole matched range with its own text collapses the\n    // three runs into one run, inheriting the formatting of the first\n    // (leading) run of the match - i.e. the Courier New / 7F6000 style.\n    found.insertText(fullText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains three paragraphs with a transcribed\n# \"<id>p074r_N</id>\" tag (N = 1, 2, 3). Each is currently split across\n# three separate runs: \"<id>\" (Courier New, color 7F6000), \"p074r_N\"\n# (plain, black) and \"</id>\" (Courier New, color 7F6000). This merges\n# each triplet back into a single run - e.g. \"<id>p074r_1</id>\" - using\n# the formatting of the leading \"<id>\" run, matching the newly\n# downloaded transcription-line (tl) formatting.\n\n$d = $word.ActiveDocument\n\nforeach ($n in 1..3) {\n    $idOpen = \"<id>\"\n    $idClose = \"</id>\"\n    $marker = \"p074r_$n\"\n    $fullText = \"$idOpen$marker$idClose\"\n\n    # Locate the full \"<id>p074r_N</id>\" span (it matches across the\n    # three run boundaries).\n    $full = $d.Content\n    $full.Find.ClearFormatting()\n    $foundFull = $full.Find.Execute($fullText)\n    if (-not $foundFull) { continue }\n\n    # Within that span, re-locate the leading \"<id>\" run - this is the\n    # run whose formatting (\"<id>\" Courier New / 7F6000) must survive.\n    $r1 = $d.Range($full.Start, $full.End)\n    $r1.Find.ClearFormatting()\n    $r1.Find.Execute($idOpen) | Out-Null\n\n    # ...and the trailing \"</id>\" run.\n    $r3 = $d.Range($r1.End, $full.End)\n    $r3.Find.ClearFormatting()\n    $r3.Find.Execute($idClose) | Out-Null\n\n    # Everything between them is the middle \"p074r_N\" run.\n    $r2 = $d.Range($r1.End, $r3.Start)\n\n    # Delete the trailing run first (so r1/r2 offsets stay valid), then\n    # the middle run, leaving only the original \"<id>\" run untouched.\n    $r3.Delete()\n    $r2.Delete()\n\n    # Append the removed text back onto the surviving \"<id>\" run so it\n    # reads \"<id>p074r_N</id>\" as a single run.\n    $r1.InsertAfter(\"$marker$idClose\")\n}\n"}
